$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5855.769
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 6030
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 18090
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -18426
$ws.Range("H74").Value = 84241.75
$ws.Range("I74").Value = 100400.414
$ws.Range("K74").Value = 100400.414
$ws.Range("M74").Value = -99464.414
$ws.Range("H77").Value = 84241.75
$ws.Range("I77").Value = 100400.414
$ws.Range("K77").Value = 502002.07
$ws.Range("M77").Value = -497322.07
$ws.Range("H99").Value = 1756.5
$ws.Range("I99").Value = 2424.75
$ws.Range("J99").Value = 420
$ws.Range("K99").Value = 7274.25
$ws.Range("L99").Value = 1260
$ws.Range("M99").Value = -5776.25
$ws.Range("N99").Value = -4256
$ws.Range("H131").Value = 7386.7
$ws.Range("I131").Value = 2675.8
$ws.Range("K131").Value = 8027.400000000001
$ws.Range("M131").Value = -2987.400000000001
$ws.Range("H138").Value = 7196.282
$ws.Range("J138").Value = 7197.237
$ws.Range("L138").Value = 21591.711
$ws.Range("N138").Value = -31871.711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1588
$ws.Range("I45").Value = 1115
$ws.Range("K45").Value = 1115
$ws.Range("M45").Value = -738
$ws.Range("H61").Value = 5735.294
$ws.Range("I61").Value = 2765.8333
$ws.Range("K61").Value = 2765.8333
$ws.Range("M61").Value = -2553.8333
$ws.Range("H74").Value = 3059.6667
$ws.Range("I74").Value = 3251.6
$ws.Range("J74").Value = 2100
$ws.Range("K74").Value = 3251.6
$ws.Range("L74").Value = 2100
$ws.Range("M74").Value = -2377.6
$ws.Range("N74").Value = -3848
$ws.Range("H77").Value = 3059.6667
$ws.Range("I77").Value = 3251.6
$ws.Range("J77").Value = 2100
$ws.Range("K77").Value = 16258
$ws.Range("L77").Value = 10500
$ws.Range("M77").Value = -11890
$ws.Range("N77").Value = -19236
$ws.Range("H122").Value = 2845.7273
$ws.Range("I122").Value = 2589.2222
$ws.Range("K122").Value = 7767.6666
$ws.Range("M122").Value = -5317.6666
$ws.Range("H133").Value = 250065
$ws.Range("J133").Value = 250065
$ws.Range("L133").Value = 250065
$ws.Range("N133").Value = -255125
$ws.Range("H136").Value = 5735.294
$ws.Range("I136").Value = 2765.8333
$ws.Range("K136").Value = 8297.499899999999
$ws.Range("M136").Value = -5747.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2967.625
$ws.Range("I20").Value = 2967.625
$ws.Range("K20").Value = 2967.625
$ws.Range("M20").Value = -2720.625
$ws.Range("H24").Value = 9500
$ws.Range("I24").Value = 9000
$ws.Range("K24").Value = 9000
$ws.Range("M24").Value = -8765
$ws.Range("H105").Value = 2168.9565
$ws.Range("I105").Value = 1737.3529
$ws.Range("K105").Value = 1737.3529
$ws.Range("M105").Value = 9.647099999999909
$ws.Range("H107").Value = 1671.1666
$ws.Range("I107").Value = 1402.8
$ws.Range("K107").Value = 1402.8
$ws.Range("M107").Value = 517.2
$ws.Range("H135").Value = 67163
$ws.Range("J135").Value = 67163
$ws.Range("L135").Value = 67163
$ws.Range("N135").Value = -77303

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12117
$ws.Range("I31").Value = 15210.7
$ws.Range("J31").Value = 8249.875
$ws.Range("K31").Value = 15210.7
$ws.Range("L31").Value = 8249.875
$ws.Range("M31").Value = -14915.7
$ws.Range("N31").Value = -8839.875
$ws.Range("H34").Value = 12117
$ws.Range("I34").Value = 15210.7
$ws.Range("J34").Value = 8249.875
$ws.Range("K34").Value = 15210.7
$ws.Range("L34").Value = 8249.875
$ws.Range("M34").Value = -15008.7
$ws.Range("N34").Value = -8653.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 28
$ws.Range("I8").Value = 28
$ws.Range("K8").Value = 84
$ws.Range("M8").Value = 55
$ws.Range("H34").Value = 4535.5557
$ws.Range("J34").Value = 4535.5557
$ws.Range("L34").Value = 13606.6671
$ws.Range("N34").Value = -13774.6671
$ws.Range("H39").Value = 3955.8
$ws.Range("J39").Value = 3955.8
$ws.Range("L39").Value = 11867.4
$ws.Range("N39").Value = -12455.4
$ws.Range("H44").Value = 10512.3
$ws.Range("I44").Value = 11647
$ws.Range("J44").Value = 300
$ws.Range("K44").Value = 34941
$ws.Range("L44").Value = 900
$ws.Range("M44").Value = -34543
$ws.Range("N44").Value = -1696
$ws.Range("H55").Value = 1769.8
$ws.Range("J55").Value = 2128.2856
$ws.Range("L55").Value = 6384.8568
$ws.Range("N55").Value = -6738.8568
$ws.Range("H136").Value = 7999.3335
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 267.82352
$ws.Range("I2").Value = 384.63635
$ws.Range("J2").Value = 53.666668
$ws.Range("K2").Value = 384.63635
$ws.Range("L2").Value = 53.666668
$ws.Range("M2").Value = -271.63635
$ws.Range("N2").Value = -279.666668
$ws.Range("H3").Value = 627132.5600000001
$ws.Range("I3").Value = 751148
$ws.Range("K3").Value = 751148
$ws.Range("M3").Value = -751032
$ws.Range("H11").Value = 540001.6
$ws.Range("I11").Value = 540001.6
$ws.Range("K11").Value = 540001.6
$ws.Range("M11").Value = -539862.6
$ws.Range("H14").Value = 3719.8
$ws.Range("I14").Value = 4266.3335
$ws.Range("J14").Value = 2900
$ws.Range("K14").Value = 4266.3335
$ws.Range("L14").Value = 2900
$ws.Range("M14").Value = -4098.3335
$ws.Range("N14").Value = -3236
$ws.Range("H36").Value = 2000
$ws.Range("J36").Value = 2000
$ws.Range("L36").Value = 2000
$ws.Range("N36").Value = -2970
$ws.Range("H70").Value = 5965.4614
$ws.Range("J70").Value = 5987.1
$ws.Range("L70").Value = 5987.1
$ws.Range("N70").Value = -6527.1
$ws.Range("H73").Value = 5965.4614
$ws.Range("J73").Value = 5987.1
$ws.Range("L73").Value = 5987.1
$ws.Range("N73").Value = -7859.1
$ws.Range("H99").Value = 8148.5
$ws.Range("I99").Value = 4178.4
$ws.Range("K99").Value = 4178.4
$ws.Range("M99").Value = -1932.4
$ws.Range("H132").Value = 2998.4
$ws.Range("I132").Value = 2937.6667
$ws.Range("K132").Value = 8813.000100000001
$ws.Range("M132").Value = -6283.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6119.6553
$ws.Range("I122").Value = 5261.8096
$ws.Range("K122").Value = 15785.4288
$ws.Range("M122").Value = -13335.4288
$ws.Range("H136").Value = 6024.567
$ws.Range("I136").Value = 5920.6816
$ws.Range("J136").Value = 6310.25
$ws.Range("K136").Value = 17762.0448
$ws.Range("L136").Value = 18930.75
$ws.Range("M136").Value = -15212.0448
$ws.Range("N136").Value = -24030.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1506.7142
$ws.Range("I81").Value = 1182.6666
$ws.Range("J81").Value = 1749.75
$ws.Range("K81").Value = 2365.3332
$ws.Range("L81").Value = 3499.5
$ws.Range("M81").Value = -1304.3332
$ws.Range("N81").Value = -5621.5
$ws.Range("H84").Value = 1506.7142
$ws.Range("I84").Value = 1182.6666
$ws.Range("J84").Value = 1749.75
$ws.Range("K84").Value = 11826.666
$ws.Range("L84").Value = 17497.5
$ws.Range("M84").Value = -6522.666000000001
$ws.Range("N84").Value = -28105.5
$ws.Range("H93").Value = 32500
$ws.Range("I93").Value = 32500
$ws.Range("K93").Value = 32500
$ws.Range("M93").Value = -30004
$ws.Range("H96").Value = 1401.5
$ws.Range("I96").Value = 1799
$ws.Range("J96").Value = 1004
$ws.Range("K96").Value = 1799
$ws.Range("L96").Value = 1004
$ws.Range("M96").Value = -426
$ws.Range("N96").Value = -3750
$ws.Range("H100").Value = 4167699.2
$ws.Range("I100").Value = 5883206
$ws.Range("J100").Value = 1468.4286
$ws.Range("K100").Value = 11766412
$ws.Range("L100").Value = 2936.8572
$ws.Range("M100").Value = -11765871
$ws.Range("N100").Value = -4018.8572
$ws.Range("H122").Value = 5500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -21400
